# Fix the typo in the dialog reader DB name: "DiglogDB" -> "DialogDB"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C6").Value = "DialogDB"

# Move/save the active selection to C13 (matches the edited view state)
$ws.Range("C13").Select()
